$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 41; $row++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value()
        if ($val -ne $null) {
            $parts = $val -split ":"
            $newparts = @()
            foreach ($p in $parts) {
                if ($p.Length -ge 2 -and $p.Substring(0,2) -eq "0x") {
                    $prefix = $p.Substring(0, 2)
                    $rest = $p.Substring(2)
                    $newparts += $prefix + $rest.ToUpper()
                } else {
                    $newparts += $p
                }
            }
            $newval = $newparts -join ":"
            $cell.Value = $newval
        }
    }
}
